# Auto-generated script applying the scheduled-runner price/profit updates
# to the Ultros_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 224.11111
$ws.Range("I2").Value = 262.85715
$ws.Range("K2").Value = 262.85715
$ws.Range("M2").Value = -149.85715

$ws.Range("H12").Value = 407.2857
$ws.Range("I12").Value = 339
$ws.Range("J12").Value = 498.33334
$ws.Range("K12").Value = 339
$ws.Range("L12").Value = 498.33334
$ws.Range("M12").Value = -169
$ws.Range("N12").Value = -838.33334

$ws.Range("H17").Value = 4808.3335
$ws.Range("J17").Value = 6462.5
$ws.Range("L17").Value = 19387.5
$ws.Range("N17").Value = -19723.5

$ws.Range("H19").Value = 1249.75
$ws.Range("J19").Value = 1078.5714
$ws.Range("L19").Value = 1078.5714
$ws.Range("N19").Value = -1428.5714

$ws.Range("H32").Value = 16671565
$ws.Range("I32").Value = 1500
$ws.Range("J32").Value = 20005578
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 20005578
$ws.Range("M32").Value = -1174
$ws.Range("N32").Value = -20006230

$ws.Range("H38").Value = 1047.8667
$ws.Range("I38").Value = 1047.8667
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 3143.6001
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -2771.6001
$ws.Range("N38").ClearContents()

$ws.Range("H40").Value = 3079.2646
$ws.Range("I40").Value = 2264.4736
$ws.Range("J40").Value = 4111.3335
$ws.Range("K40").Value = 2264.4736
$ws.Range("L40").Value = 4111.3335
$ws.Range("M40").Value = -2089.4736
$ws.Range("N40").Value = -4461.3335

$ws.Range("H62").Value = 1993.5
$ws.Range("J62").Value = 1993.5
$ws.Range("L62").Value = 1993.5
$ws.Range("N62").Value = -3241.5

$ws.Range("H65").Value = 1993.5
$ws.Range("J65").Value = 1993.5
$ws.Range("L65").Value = 9967.5
$ws.Range("N65").Value = -16207.5

$ws.Range("H74").Value = 14333.333
$ws.Range("J74").Value = 20000
$ws.Range("L74").Value = 20000
$ws.Range("N74").Value = -21872

$ws.Range("H77").Value = 14333.333
$ws.Range("J77").Value = 20000
$ws.Range("L77").Value = 100000
$ws.Range("N77").Value = -109360

$ws.Range("H86").Value = 2291307
$ws.Range("I86").Value = 2785.6428
$ws.Range("J86").Value = 5851229.5
$ws.Range("K86").Value = 2785.6428
$ws.Range("L86").Value = 5851229.5
$ws.Range("M86").Value = -1662.6428
$ws.Range("N86").Value = -5853475.5

$ws.Range("H87").Value = 33168.25
$ws.Range("J87").Value = 33168.25
$ws.Range("L87").Value = 33168.25
$ws.Range("N87").Value = -35664.25

$ws.Range("H89").Value = 2291307
$ws.Range("I89").Value = 2785.6428
$ws.Range("J89").Value = 5851229.5
$ws.Range("K89").Value = 13928.214
$ws.Range("L89").Value = 29256147.5
$ws.Range("M89").Value = -8312.214
$ws.Range("N89").Value = -29267379.5

$ws.Range("H90").Value = 33168.25
$ws.Range("J90").Value = 33168.25
$ws.Range("L90").Value = 99504.75
$ws.Range("N90").Value = -111984.75

$ws.Range("H93").Value = 37500
$ws.Range("J93").Value = 37500
$ws.Range("L93").Value = 37500
$ws.Range("N93").Value = -42492

$ws.Range("H98").Value = 3000
$ws.Range("I98").Value = 3000
$ws.Range("K98").Value = 3000
$ws.Range("M98").Value = -1502

$ws.Range("H103").Value = 2414.1667
$ws.Range("I103").Value = 1745
$ws.Range("J103").Value = 2748.75
$ws.Range("K103").Value = 5235
$ws.Range("L103").Value = 8246.25
$ws.Range("M103").Value = -4649
$ws.Range("N103").Value = -9418.25

$ws.Range("H106").Value = 2812.25
$ws.Range("I106").Value = 2985
$ws.Range("J106").Value = 1603
$ws.Range("K106").Value = 2985
$ws.Range("L106").Value = 1603
$ws.Range("M106").Value = -2354
$ws.Range("N106").Value = -2865

$ws.Range("H121").Value = 3999.4285
$ws.Range("J121").Value = 3999.4285
$ws.Range("L121").Value = 11998.2855
$ws.Range("N121").Value = -15492.2855

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550

$ws.Range("H123").Value = 59999.59
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

$ws.Range("H125").Value = 15992.143
$ws.Range("I125").Value = 2750
$ws.Range("J125").Value = 18199.166
$ws.Range("K125").Value = 24750
$ws.Range("L125").Value = 163792.494
$ws.Range("M125").Value = -22290
$ws.Range("N125").Value = -168712.494

$ws.Range("H129").Value = 135257.94
$ws.Range("I129").Value = 144562.08
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 433686.24
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = -428686.24
$ws.Range("N129").Value = -25000

$ws.Range("H132").Value = 47351.727
$ws.Range("I132").Value = 2273.7778
$ws.Range("J132").Value = 250202.5
$ws.Range("K132").Value = 6821.3334
$ws.Range("L132").Value = 750607.5
$ws.Range("M132").Value = -4291.3334
$ws.Range("N132").Value = -755667.5

$ws.Range("H135").Value = 2980
$ws.Range("I135").Value = 2980
$ws.Range("K135").Value = 26820
$ws.Range("M135").Value = -24285

$ws.Range("H137").Value = 5702.3477
$ws.Range("I137").Value = 7867.3076
$ws.Range("J137").Value = 2887.9
$ws.Range("K137").Value = 23601.9228
$ws.Range("L137").Value = 8663.7
$ws.Range("M137").Value = -21051.9228
$ws.Range("N137").Value = -13763.7

$ws.Range("H138").Value = 2938
$ws.Range("I138").Value = 2133.1667
$ws.Range("J138").Value = 3541.625
$ws.Range("K138").Value = 6399.500100000001
$ws.Range("L138").Value = 10624.875
$ws.Range("M138").Value = -1259.500100000001
$ws.Range("N138").Value = -20904.875

$ws.Range("H141").Value = 4679.1875
$ws.Range("I141").Value = 3324.4666
$ws.Range("K141").Value = 9973.399800000001
$ws.Range("M141").Value = -4793.399800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 2500
$ws.Range("I6").Value = 2500
$ws.Range("K6").Value = 2500
$ws.Range("M6").Value = -2327

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

$ws.Range("H26").Value = 1153.1666
$ws.Range("I26").Value = 783.8
$ws.Range("J26").Value = 3000
$ws.Range("K26").Value = 783.8
$ws.Range("L26").Value = 3000
$ws.Range("M26").Value = -453.8
$ws.Range("N26").Value = -3660

$ws.Range("H41").Value = 21666.334
$ws.Range("I41").Value = 21666.334
$ws.Range("K41").Value = 21666.334
$ws.Range("M41").Value = -21252.334

$ws.Range("H45").Value = 1704
$ws.Range("I45").Value = 1153.8334
$ws.Range("J45").Value = 2804.3333
$ws.Range("K45").Value = 1153.8334
$ws.Range("L45").Value = 2804.3333
$ws.Range("M45").Value = -776.8334
$ws.Range("N45").Value = -3558.3333

$ws.Range("H61").Value = 4463.9165
$ws.Range("I61").Value = 4248.923
$ws.Range("J61").Value = 4718
$ws.Range("K61").Value = 4248.923
$ws.Range("L61").Value = 4718
$ws.Range("M61").Value = -4036.923
$ws.Range("N61").Value = -5142

$ws.Range("H74").Value = 2092.3333
$ws.Range("I74").Value = 2191.2
$ws.Range("J74").Value = 1598
$ws.Range("K74").Value = 2191.2
$ws.Range("L74").Value = 1598
$ws.Range("M74").Value = -1317.2
$ws.Range("N74").Value = -3346

$ws.Range("H77").Value = 2092.3333
$ws.Range("I77").Value = 2191.2
$ws.Range("J77").Value = 1598
$ws.Range("K77").Value = 10956
$ws.Range("L77").Value = 7990
$ws.Range("M77").Value = -6588
$ws.Range("N77").Value = -16726

$ws.Range("H88").Value = 7577102.5
$ws.Range("I88").Value = 1144.75
$ws.Range("J88").Value = 11906221
$ws.Range("K88").Value = 1144.75
$ws.Range("L88").Value = 11906221
$ws.Range("M88").Value = -738.75
$ws.Range("N88").Value = -11907033

$ws.Range("H91").Value = 7577102.5
$ws.Range("I91").Value = 1144.75
$ws.Range("J91").Value = 11906221
$ws.Range("K91").Value = 1144.75
$ws.Range("L91").Value = 11906221
$ws.Range("M91").Value = 259.25
$ws.Range("N91").Value = -11909029

$ws.Range("H97").Value = 4274378
$ws.Range("J97").Value = 13889369
$ws.Range("L97").Value = 13889369
$ws.Range("N97").Value = -13890361

$ws.Range("H102").Value = 17545858
$ws.Range("I102").Value = 2188.2856
$ws.Range("K102").Value = 2188.2856
$ws.Range("M102").Value = -566.2856000000002

$ws.Range("H110").Value = 840.8182
$ws.Range("I110").Value = 834.9
$ws.Range("K110").Value = 834.9
$ws.Range("M110").Value = 1210.1

$ws.Range("H122").Value = 3574.0715
$ws.Range("I122").Value = 1962.4736
$ws.Range("J122").Value = 6976.3335
$ws.Range("K122").Value = 5887.4208
$ws.Range("L122").Value = 20929.0005
$ws.Range("M122").Value = -3437.4208
$ws.Range("N122").Value = -25829.0005

$ws.Range("H131").Value = 70000
$ws.Range("J131").Value = 70000
$ws.Range("L131").Value = 70000
$ws.Range("N131").Value = -80080

$ws.Range("H132").Value = 3077.5186
$ws.Range("I132").Value = 1898.5333
$ws.Range("J132").Value = 4551.25
$ws.Range("K132").Value = 5695.5999
$ws.Range("L132").Value = 13653.75
$ws.Range("M132").Value = -3165.5999
$ws.Range("N132").Value = -18713.75

$ws.Range("H136").Value = 4463.9165
$ws.Range("I136").Value = 4248.923
$ws.Range("J136").Value = 4718
$ws.Range("K136").Value = 12746.769
$ws.Range("L136").Value = 14154
$ws.Range("M136").Value = -10196.769
$ws.Range("N136").Value = -19254

$ws.Range("H137").Value = 80000
$ws.Range("J137").Value = 80000
$ws.Range("L137").Value = 80000
$ws.Range("N137").Value = -90200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1666798.4
$ws.Range("I7").Value = 197.5
$ws.Range("J7").Value = 5000000
$ws.Range("K7").Value = 197.5
$ws.Range("L7").Value = 5000000
$ws.Range("M7").Value = -84.5
$ws.Range("N7").Value = -5000226

$ws.Range("H11").Value = 592.375
$ws.Range("J11").Value = 1131.75
$ws.Range("L11").Value = 1131.75
$ws.Range("N11").Value = -1411.75

$ws.Range("H22").Value = 2686.818
$ws.Range("I22").Value = 4199.3335
$ws.Range("J22").Value = 871.8
$ws.Range("K22").Value = 4199.3335
$ws.Range("L22").Value = 871.8
$ws.Range("M22").Value = -4026.3335
$ws.Range("N22").Value = -1217.8

$ws.Range("H23").Value = 25000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H29").Value = 5140.3335
$ws.Range("I29").Value = 210.5
$ws.Range("K29").Value = 210.5
$ws.Range("M29").Value = 78.5

$ws.Range("H36").Value = 3000
$ws.Range("I36").Value = 2000
$ws.Range("K36").Value = 2000
$ws.Range("M36").Value = -1466

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H86").Value = 5417.636
$ws.Range("I86").Value = 4719.4
$ws.Range("J86").Value = 5999.5
$ws.Range("K86").Value = 4719.4
$ws.Range("L86").Value = 5999.5
$ws.Range("M86").Value = -3596.4
$ws.Range("N86").Value = -8245.5

$ws.Range("H89").Value = 5417.636
$ws.Range("I89").Value = 4719.4
$ws.Range("J89").Value = 5999.5
$ws.Range("K89").Value = 23597
$ws.Range("L89").Value = 29997.5
$ws.Range("M89").Value = -17981
$ws.Range("N89").Value = -41229.5

$ws.Range("H99").Value = 1901.9395
$ws.Range("I99").Value = 1908.4
$ws.Range("J99").Value = 1892
$ws.Range("K99").Value = 1908.4
$ws.Range("L99").Value = 1892
$ws.Range("M99").Value = -410.4000000000001
$ws.Range("N99").Value = -4888

$ws.Range("H105").Value = 3869.6191
$ws.Range("I105").Value = 3393.6428
$ws.Range("J105").Value = 4821.5713
$ws.Range("K105").Value = 3393.6428
$ws.Range("L105").Value = 4821.5713
$ws.Range("M105").Value = -1646.6428
$ws.Range("N105").Value = -8315.5713

$ws.Range("H107").Value = 1803.5
$ws.Range("I107").Value = 1811.1428
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 1811.1428
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = 108.8571999999999
$ws.Range("N107").Value = -5590

$ws.Range("H134").Value = 4139.3945
$ws.Range("I134").Value = 3029.8076
$ws.Range("J134").Value = 6543.5
$ws.Range("K134").Value = 9089.4228
$ws.Range("L134").Value = 19630.5
$ws.Range("M134").Value = -6554.4228
$ws.Range("N134").Value = -24700.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3500342.8
$ws.Range("I6").Value = 3500342.8
$ws.Range("K6").Value = 3500342.8
$ws.Range("M6").Value = -3500229.8

$ws.Range("H12").Value = 3800
$ws.Range("J12").Value = 5500
$ws.Range("L12").Value = 5500
$ws.Range("N12").Value = -5840

$ws.Range("H31").Value = 3186.25
$ws.Range("I31").Value = 2179.2144
$ws.Range("J31").Value = 3827.0908
$ws.Range("K31").Value = 2179.2144
$ws.Range("L31").Value = 3827.0908
$ws.Range("M31").Value = -1884.2144
$ws.Range("N31").Value = -4417.0908

$ws.Range("H34").Value = 3186.25
$ws.Range("I34").Value = 2179.2144
$ws.Range("J34").Value = 3827.0908
$ws.Range("K34").Value = 2179.2144
$ws.Range("L34").Value = 3827.0908
$ws.Range("M34").Value = -1977.2144
$ws.Range("N34").Value = -4231.0908

$ws.Range("H53").Value = 44993.332
$ws.Range("J53").Value = 44993.332
$ws.Range("L53").Value = 44993.332
$ws.Range("N53").Value = -46207.332

$ws.Range("H74").Value = 53815.4
$ws.Range("J74").Value = 53815.4
$ws.Range("L74").Value = 53815.4
$ws.Range("N74").Value = -55563.4

$ws.Range("H77").Value = 53815.4
$ws.Range("J77").Value = 53815.4
$ws.Range("L77").Value = 161446.2
$ws.Range("N77").Value = -170182.2

$ws.Range("H86").Value = 29043.5
$ws.Range("I86").Value = 54365
$ws.Range("J86").Value = 10052.375
$ws.Range("K86").Value = 54365
$ws.Range("L86").Value = 10052.375
$ws.Range("M86").Value = -53242
$ws.Range("N86").Value = -12298.375

$ws.Range("H89").Value = 29043.5
$ws.Range("I89").Value = 54365
$ws.Range("J89").Value = 10052.375
$ws.Range("K89").Value = 271825
$ws.Range("L89").Value = 50261.875
$ws.Range("M89").Value = -266209
$ws.Range("N89").Value = -61493.875

$ws.Range("H99").Value = 9355415
$ws.Range("J99").Value = 18189228
$ws.Range("L99").Value = 18189228
$ws.Range("N99").Value = -18192224

$ws.Range("H105").Value = 6513.4375
$ws.Range("I105").Value = 1179.7778
$ws.Range("J105").Value = 13371
$ws.Range("K105").Value = 1179.7778
$ws.Range("L105").Value = 13371
$ws.Range("M105").Value = 567.2221999999999
$ws.Range("N105").Value = -16865

$ws.Range("H110").Value = 30702
$ws.Range("J110").Value = 30702
$ws.Range("L110").Value = 30702
$ws.Range("N110").Value = -38882

$ws.Range("H112").Value = 39999.8
$ws.Range("J112").Value = 39999.8
$ws.Range("L112").Value = 39999.8
$ws.Range("N112").Value = -42953.8

$ws.Range("H114").Value = 65000
$ws.Range("J114").Value = 65000
$ws.Range("L114").Value = 65000
$ws.Range("N114").Value = -73678

$ws.Range("H115").Value = 39999.832
$ws.Range("J115").Value = 39999.832
$ws.Range("L115").Value = 39999.832
$ws.Range("N115").Value = -42349.832

$ws.Range("H120").Value = 892999.8
$ws.Range("J120").Value = 892999.8
$ws.Range("L120").Value = 892999.8
$ws.Range("N120").Value = -900257.8

$ws.Range("H126").Value = 9355415
$ws.Range("J126").Value = 18189228
$ws.Range("L126").Value = 54567684
$ws.Range("N126").Value = -54572624

$ws.Range("H133").Value = 62259.125
$ws.Range("J133").Value = 66867.57
$ws.Range("L133").Value = 66867.57
$ws.Range("N133").Value = -71927.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 300.2
$ws.Range("J2").Value = 400.33334
$ws.Range("L2").Value = 2402.00004
$ws.Range("N2").Value = -2628.00004

$ws.Range("H7").Value = 1081.8667
$ws.Range("I7").Value = 1285.5
$ws.Range("J7").Value = 267.33334
$ws.Range("K7").Value = 3856.5
$ws.Range("L7").Value = 802.0000200000001
$ws.Range("M7").Value = -3744.5
$ws.Range("N7").Value = -1026.00002

$ws.Range("H11").Value = 328.8
$ws.Range("I11").Value = 439.66666
$ws.Range("J11").Value = 162.5
$ws.Range("K11").Value = 1318.99998
$ws.Range("L11").Value = 487.5
$ws.Range("M11").Value = -1178.99998
$ws.Range("N11").Value = -767.5

$ws.Range("H26").Value = 122.181816
$ws.Range("J26").Value = 228
$ws.Range("L26").Value = 684
$ws.Range("N26").Value = -1260

$ws.Range("H80").Value = 2807.8
$ws.Range("J80").Value = 2919.0715
$ws.Range("L80").Value = 8757.2145
$ws.Range("N80").Value = -10629.2145

$ws.Range("H83").Value = 2807.8
$ws.Range("J83").Value = 2919.0715
$ws.Range("L83").Value = 26271.6435
$ws.Range("N83").Value = -35631.6435

$ws.Range("H113").Value = 2226.2
$ws.Range("J113").Value = 2367.6667
$ws.Range("L113").Value = 7103.000100000001
$ws.Range("N113").Value = -11443.0001

$ws.Range("H114").Value = 2257
$ws.Range("I114").Value = 900
$ws.Range("J114").Value = 2483.1667
$ws.Range("K114").Value = 2700
$ws.Range("L114").Value = 7449.500100000001
$ws.Range("M114").Value = 554
$ws.Range("N114").Value = -13957.5001

$ws.Range("H117").Value = 2049.6667
$ws.Range("I117").Value = 649
$ws.Range("K117").Value = 1947
$ws.Range("M117").Value = 1495

$ws.Range("H121").Value = 1259.5555
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 1639.3334
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 4918.0002
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -7538.0002

$ws.Range("H131").Value = 4411.25
$ws.Range("I131").Value = 2558.6667
$ws.Range("J131").Value = 5927
$ws.Range("K131").Value = 7676.000100000001
$ws.Range("L131").Value = 17781
$ws.Range("M131").Value = -2636.000100000001
$ws.Range("N131").Value = -27861

$ws.Range("H134").Value = 6922.7
$ws.Range("I134").Value = 1461
$ws.Range("K134").Value = 4383
$ws.Range("M134").Value = 687

$ws.Range("H137").Value = 3356.1667
$ws.Range("I137").Value = 2165
$ws.Range("J137").Value = 5738.5
$ws.Range("K137").Value = 6495
$ws.Range("L137").Value = 17215.5
$ws.Range("M137").Value = -1395
$ws.Range("N137").Value = -27415.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2373466.8
$ws.Range("J14").Value = 1980171.4
$ws.Range("L14").Value = 1980171.4
$ws.Range("N14").Value = -1980507.4

$ws.Range("H33").Value = 24974
$ws.Range("I33").Value = 24973
$ws.Range("J33").Value = 24974.334
$ws.Range("K33").Value = 24973
$ws.Range("L33").Value = 24974.334
$ws.Range("M33").Value = -24721
$ws.Range("N33").Value = -25478.334

$ws.Range("H36").Value = 975
$ws.Range("I36").Value = 975
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 975
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -490
$ws.Range("N36").ClearContents()

$ws.Range("H53").Value = 15000
$ws.Range("J53").Value = 15000
$ws.Range("L53").Value = 15000
$ws.Range("N53").Value = -16262

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H97").Value = 16666850
$ws.Range("I97").Value = 200.54546
$ws.Range("K97").Value = 200.54546
$ws.Range("M97").Value = 295.45454

$ws.Range("H107").Value = 358.77777
$ws.Range("J107").Value = 458
$ws.Range("L107").Value = 458
$ws.Range("N107").Value = -4298

$ws.Range("H132").Value = 2867.8333
$ws.Range("I132").Value = 735.6667
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 2207.0001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = 322.9998999999998
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 324.5
$ws.Range("J12").Value = 249
$ws.Range("L12").Value = 249
$ws.Range("N12").Value = -589

$ws.Range("H16").Value = 179.07143
$ws.Range("I16").Value = 192.25
$ws.Range("K16").Value = 192.25
$ws.Range("M16").Value = -22.25

$ws.Range("H40").Value = 6496.3823
$ws.Range("I40").Value = 8017.75
$ws.Range("J40").Value = 5144.0557
$ws.Range("K40").Value = 8017.75
$ws.Range("L40").Value = 5144.0557
$ws.Range("M40").Value = -7881.75
$ws.Range("N40").Value = -5416.0557

$ws.Range("H43").Value = 10000000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H55").Value = 1765.3572
$ws.Range("I55").Value = 3118.1428
$ws.Range("K55").Value = 3118.1428
$ws.Range("M55").Value = -2945.1428

$ws.Range("H61").Value = 4182.8887
$ws.Range("I61").Value = 3687
$ws.Range("K61").Value = 3687
$ws.Range("M61").Value = -3485

$ws.Range("H68").Value = 220708.3
$ws.Range("J68").Value = 380036.75
$ws.Range("L68").Value = 380036.75
$ws.Range("N68").Value = -381534.75

$ws.Range("H71").Value = 220708.3
$ws.Range("J71").Value = 380036.75
$ws.Range("L71").Value = 1900183.75
$ws.Range("N71").Value = -1907671.75

$ws.Range("H100").Value = 4137.579
$ws.Range("I100").Value = 2151.25
$ws.Range("J100").Value = 5582.1816
$ws.Range("K100").Value = 2151.25
$ws.Range("L100").Value = 5582.1816
$ws.Range("M100").Value = -1610.25
$ws.Range("N100").Value = -6664.1816

$ws.Range("H113").Value = 4182.8887
$ws.Range("I113").Value = 3687
$ws.Range("K113").Value = 3687
$ws.Range("M113").Value = -1517

$ws.Range("H122").Value = 5542.154
$ws.Range("I122").Value = 2560.2222
$ws.Range("J122").Value = 7120.8237
$ws.Range("K122").Value = 7680.6666
$ws.Range("L122").Value = 21362.4711
$ws.Range("M122").Value = -5230.6666
$ws.Range("N122").Value = -26262.4711

$ws.Range("H132").Value = 4847.884
$ws.Range("I132").Value = 3240.7812
$ws.Range("J132").Value = 9523.091
$ws.Range("K132").Value = 9722.3436
$ws.Range("L132").Value = 28569.273
$ws.Range("M132").Value = -7192.3436
$ws.Range("N132").Value = -33629.273

$ws.Range("H136").Value = 6746.9287
$ws.Range("I136").Value = 4814.5
$ws.Range("J136").Value = 7519.9
$ws.Range("K136").Value = 14443.5
$ws.Range("L136").Value = 22559.7
$ws.Range("M136").Value = -11893.5
$ws.Range("N136").Value = -27659.7

$ws.Range("H137").Value = 69999.52
$ws.Range("J137").Value = 69999.52
$ws.Range("L137").Value = 69999.52
$ws.Range("N137").Value = -80199.52

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 35591
$ws.Range("I47").Value = 35464.066
$ws.Range("J47").Value = 37495
$ws.Range("K47").Value = 35464.066
$ws.Range("L47").Value = 37495
$ws.Range("M47").Value = -34892.066
$ws.Range("N47").Value = -38639

$ws.Range("H68").Value = 200271
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 200271
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 200271
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -201893

$ws.Range("H71").Value = 200271
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 200271
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 600813
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -608925

$ws.Range("H122").Value = 266622.7
$ws.Range("I122").Value = 387151.3
$ws.Range("K122").Value = 1161453.9
$ws.Range("M122").Value = -1159003.9

$ws.Range("H126").Value = 2411
$ws.Range("I126").Value = 2411
$ws.Range("K126").Value = 7233
$ws.Range("M126").Value = -4763

$ws.Range("H132").Value = 1839.2778
$ws.Range("I132").Value = 945.375
$ws.Range("J132").Value = 3627.0833
$ws.Range("K132").Value = 2836.125
$ws.Range("L132").Value = 10881.2499
$ws.Range("M132").Value = -306.125
$ws.Range("N132").Value = -15941.2499

$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
